$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format for Price cells whose new values would otherwise
# be auto-detected as numbers by Excel (so they stay text, matching the source data).
$textCells = @("D6","D7","D8","D10","D11","D12","D13","D14","D15","D18","D19","D21","D22","D23","D24","D26","D27","D28","D29","D30","D31","D33","D34","D36","D37","D38","D39","D40","D41","D42","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = "22.395.46"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.574.07"
$ws.Range("E3").Value = "  +0.17%  "
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").Value = "291.22"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").Value = "0.3767"
$ws.Range("E7").Value = "  +3.03%  "
$ws.Range("D8").Value = "50.18"
$ws.Range("E8").Value = "  +1.49%  "
$ws.Range("E9").Value = "  +3.05%  "
$ws.Range("D10").Value = "1.167"
$ws.Range("E10").Value = "  +0.51%  "
$ws.Range("D11").Value = "0.07688"
$ws.Range("E11").Value = "  +1.98%  "
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("D13").Value = "21.41"
$ws.Range("E13").Value = "  +2.28%  "
$ws.Range("D14").Value = "5.986"
$ws.Range("E14").Value = "  -1.43%  "
$ws.Range("D15").Value = "6.930"
$ws.Range("E15").Value = "  +1.44%  "
$ws.Range("E16").Value = "  +1.26%  "
$ws.Range("D17").Value = "1.574.46"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "90.42"
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("D19").Value = "0.06732"
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("D21").Value = "16.81"
$ws.Range("E21").Value = "  +3.50%  "
$ws.Range("D22").Value = "6.243"
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("D23").Value = "0.5287"
$ws.Range("E23").Value = "  -3.95%  "
$ws.Range("D24").Value = "12.02"
$ws.Range("E24").Value = "  +1.19%  "
$ws.Range("D25").Value = "22.403.94"
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").Value = "2.394"
$ws.Range("E26").Value = "  +1.19%  "
$ws.Range("D27").Value = "2.769"
$ws.Range("E27").Value = "  -3.71%  "
$ws.Range("D28").Value = "20.30"
$ws.Range("E28").Value = "  +3.14%  "
$ws.Range("D29").Value = "144.75"
$ws.Range("E29").Value = "  -0.56%  "
$ws.Range("D30").Value = "5.088"
$ws.Range("E30").Value = "  +2.85%  "
$ws.Range("D31").Value = "126.43"
$ws.Range("E31").Value = "  +1.32%  "
$ws.Range("D32").Value = "1.749.29"
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "6.256"
$ws.Range("E33").Value = "  +0.85%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "1.022"
$ws.Range("E34").Value = "  +6.80%  "
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("D36").Value = "10.09"
$ws.Range("E36").Value = "  -1.81%  "
$ws.Range("D37").Value = "0.08550"
$ws.Range("E37").Value = "  +0.39%  "
$ws.Range("D38").Value = "0.02558"
$ws.Range("E38").Value = "  +2.43%  "
$ws.Range("D39").Value = "0.2329"
$ws.Range("E39").Value = "  +2.08%  "
$ws.Range("D40").Value = "0.06566"
$ws.Range("E40").Value = "  +0.88%  "
$ws.Range("D41").Value = "5.515"
$ws.Range("E41").Value = "  +1.67%  "
$ws.Range("D42").Value = "1.295"
$ws.Range("E42").Value = "  +2.36%  "
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("D44").Value = "0.6469"
$ws.Range("E44").Value = "  +2.23%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "14.19"
$ws.Range("E45").Value = "  -2.01%  "
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").Value = "1.002"
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("D47").Value = "0.6034"
$ws.Range("E47").Value = "  +1.57%  "
$ws.Range("D48").Value = "3.782"
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("D49").Value = "1.304"
$ws.Range("E49").Value = "  +11.61%  "
$ws.Range("D50").Value = "2.104"
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("D51").Value = "125.29"
$ws.Range("E51").Value = "  +3.08%  "
